$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (prices) are written as text, matching the
# original inline-string cell contents (avoids Excel auto-converting to numbers).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.001.35'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.318.04'
$ws.Range('E3').Value = '  +2.63%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '253.73'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.642'
$ws.Range('E6').Value = '  +0.65%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '76.13'
$ws.Range('E7').Value = '  +7.25%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.651'
$ws.Range('E9').Value = '  -2.55%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.68'
$ws.Range('E10').Value = '  +0.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0990'
$ws.Range('E11').Value = '  +1.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.58'
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('E13').Value = '  +1.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.665.96'
$ws.Range('E14').Value = '  +2.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.48'
$ws.Range('E15').Value = '  +4.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.886'
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.322.20'
$ws.Range('E17').Value = '  +2.58%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.988.97'
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('E19').Value = '  +3.15%  '
$ws.Range('E20').Value = '  +1.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.98'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '236.87'
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('E23').Value = '  +5.56%  '
$ws.Range('E24').Value = '  -0.69%  '
$ws.Range('E25').Value = '  -0.56%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').Value = '  -0.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.26'
$ws.Range('E28').Value = '  +1.83%  '
$ws.Range('E29').Value = '  +0.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '167.40'
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.29'
$ws.Range('E31').Value = '  +0.78%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0846'
$ws.Range('E32').Value = '  +9.65%  '
$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.128'
$ws.Range('E34').Value = '  +2.10%  '
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '30.39'
$ws.Range('E35').Value = '  +5.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.60'
$ws.Range('E36').Value = '  +11.47%  '
$ws.Range('E37').Value = '  +3.56%  '
$ws.Range('E38').Value = '  -1.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.07'
$ws.Range('E39').Value = '  +16.02%  '
$ws.Range('E40').Value = '  +2.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.90'
$ws.Range('E41').Value = '  +1.32%  '
$ws.Range('E42').Value = '  +8.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '62.72'
$ws.Range('E43').Value = '  -2.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '9.22'
$ws.Range('E44').Value = '  +3.50%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.89'
$ws.Range('E45').Value = '  -2.82%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '107.34'
$ws.Range('E46').Value = '  +13.16%  '
$ws.Range('E47').Value = '  -0.22%  '
$ws.Range('E48').Value = '  +0.76%  '
$ws.Range('E49').Value = '  -0.31%  '
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.33'
$ws.Range('E51').Value = '  -0.86%  '
